$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.247.85"
$ws.Range("E2").Value = "  +2.42%  "

$ws.Range("D3").Value = "1.880.25"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "316.58"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "0.4323"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("D8").Value = "0.3715"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07427"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "0.8898"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "21.24"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.929.70"
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.496"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.658"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.06996"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "81.38"
$ws.Range("E17").Value = "  +3.26%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000009163"
$ws.Range("E18").Value = "  +3.07%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "15.68"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "28.310.99"
$ws.Range("E21").Value = "  +2.57%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.103"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "11.06"
$ws.Range("E23").Value = "  +3.60%  "

$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.132.24"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.990"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "154.54"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "5.470"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "118.45"
$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "1.907"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08986"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.8013"
$ws.Range("E32").Value = "  +5.07%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.712"
$ws.Range("E33").Value = "  +3.00%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.181"
$ws.Range("E34").Value = "  +7.39%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.975"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.135"
$ws.Range("E36").Value = "  +3.51%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05489"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01965"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.910"
$ws.Range("E40").Value = "  +3.65%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.1703"
$ws.Range("E41").Value = "  +2.80%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5186"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.904"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "8.628"
$ws.Range("E44").Value = "  +4.29%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.06619"
$ws.Range("E46").Value = "  +0.59%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4788"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "105.91"
$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.667"
$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.859"
$ws.Range("E51").Value = "  +6.77%  "
